$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the date column (B) for new rows is treated as text, matching existing data
$ws.Range("B36:B54").NumberFormat = "@"

$data = @(
    @(14731577, "2025-09-22", "Benoit Paire", "Edas Butvilas", "Gana Benoit Paire", 4),
    @(14654903, "2025-09-21", "Jay Clarke", "Franco Agamenone", "Gana Jay Clarke", 2.2),
    @(14730387, "2025-09-21", "Jakub Paul", "Maxime Janvier", "Gana Maxime Janvier", 2.62),
    @(14730385, "2025-09-21", "Cyril Vandermeersch", "Norbert Gombos", "Gana Cyril Vandermeersch", 2.75),
    @(14730383, "2025-09-21", "Patrick Zahraj", "Loann Massard", "Gana Loann Massard", 3.75),
    @(14732347, "2025-09-21", "Maximus Jones", "Fernando Cavallo", "Gana Fernando Cavallo", 3.5),
    @(14732343, "2025-09-21", "Gonzalo Villanueva", "Santiago De la Fuente", "Gana Santiago De la Fuente", 2.62),
    @(14731191, "2025-09-21", "Maxime Chazal", "Diego Dedura-Palomero", "Gana Maxime Chazal", 2.75),
    @(14731336, "2025-09-21", "Sumit Nagal", "Alexander Donski", "Gana Alexander Donski", 3.5),
    @(14731337, "2025-09-21", "Tiago Torres", "Christoph Negritu", "Gana Tiago Torres", 3.25),
    @(14733387, "2025-09-21", "Darian King", "Tristan McCormick", "Gana Darian King", 2),
    @(14733408, "2025-09-21", "Ryuki Matsuda", "Erik Arutiunian", "Gana Erik Arutiunian", 1.57),
    @(14733401, "2025-09-21", "Carlo Alberto Caniato", "Strong Kirchheimer", "Gana Strong Kirchheimer", 2.62),
    @(14733402, "2025-09-21", "Maxwell McKennon", "Joshua Sheehy", "Gana Joshua Sheehy", 1.83),
    @(14733409, "2025-09-21", "Quinn Vandecasteele", "Bor Artnak", "Gana Quinn Vandecasteele", 1.83),
    @(14733407, "2025-09-21", "Daniel Milavsky", "Luca Castagnola", "Gana Luca Castagnola", 2.75),
    @(14733405, "2025-09-21", "Renta Tokuda", "Ramkumar Ramanathan", "Gana Ramkumar Ramanathan", 3.25),
    @(14733258, "2025-09-21", "Toby Samuel", "Liam Broady", "Gana Toby Samuel", 2.2),
    @(14731508, "2025-09-21", "Daniil Glinka", "Maxence Beauge", "Gana Maxence Beauge", 3.25)
)

$r = 36
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Range($ws.Cells.Item($r, 7), $ws.Cells.Item($r, 8)).Value = ""
    $ws.Range($ws.Cells.Item($r, 7), $ws.Cells.Item($r, 8)).ClearFormats()
    $r++
}